# result_analysis.xlsx update
# - Removed "not solved" print-style messages, replaced with a more
#   descriptive "not solved within 10min timeframe" message.
# - Highlighted the optimal "Plan length" result cells with a green fill.
# - Appended the results of four more searches (problems 6-10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Highlight the "optimal" plan-length cells (green fill) for the
#    existing result tables.
# ---------------------------------------------------------------------
$greenCells = @("E3","E4","E5","E9","E27","E28","E29")
foreach ($addr in $greenCells) {
    $ws.Range($addr).Interior.Color = 5296274   # RGB(146,208,80) = FF92D050
}

# ---------------------------------------------------------------------
# 2. Replace the "not solved" messages with a clearer message.
# ---------------------------------------------------------------------
$notSolvedCells = @("B10","B11","B22","B23")
foreach ($addr in $notSolvedCells) {
    $ws.Range($addr).Value = "not solved within 10min timeframe"
}

# ---------------------------------------------------------------------
# 3. Append the new search results (problems 6 - 10).
# ---------------------------------------------------------------------

function Set-Header($row) {
    $ws.Cells.Item($row, 1).Value = "Problem"
    $ws.Cells.Item($row, 2).Value = "Expansions"
    $ws.Cells.Item($row, 3).Value = "Goal Tests"
    $ws.Cells.Item($row, 4).Value = "New Nodes"
    $ws.Cells.Item($row, 5).Value = "Plan length"
    $ws.Cells.Item($row, 6).Value = "Time"
}

# --- 6. Recursive best first search H1 --------------------------------
$ws.Range("A31").Value = "6. Recursive best first search H1"
Set-Header 32

$ws.Cells.Item(33, 1).Value = 1
$ws.Cells.Item(33, 2).Value = 4229
$ws.Cells.Item(33, 3).Value = 4330
$ws.Cells.Item(33, 4).Value = 17023
$ws.Cells.Item(33, 5).Value = 6
$ws.Range("E33").Interior.Color = 5296274
$ws.Cells.Item(33, 6).Value = 3.42

$ws.Cells.Item(34, 1).Value = 2
$ws.Cells.Item(35, 1).Value = 3

# --- 7. --------------------------------------------------------------
# "7." alone would be auto-coerced to the number 7 by the Value setter,
# so the cell is briefly marked as Text to force a string, then restored
# to the workbook's normal style (matching the author's plain text cell).
$ws.Range("A37").NumberFormat = "@"
$ws.Range("A37").Value = "7."
$ws.Range("A37").Style = "Stand."
Set-Header 38

$ws.Cells.Item(39, 1).Value = 1
$ws.Cells.Item(40, 1).Value = 2
$ws.Cells.Item(41, 1).Value = 3

# --- 8. --------------------------------------------------------------
$ws.Range("A43").NumberFormat = "@"
$ws.Range("A43").Value = "8."
$ws.Range("A43").Style = "Stand."
Set-Header 44

$ws.Cells.Item(45, 1).Value = 1
$ws.Cells.Item(46, 1).Value = 2
$ws.Cells.Item(47, 1).Value = 3

# --- 9. astar_search h_ignore_preconditions ---------------------------
$ws.Range("A49").Value = "9. astar_search h_ignore_preconditions"
Set-Header 50

$ws.Cells.Item(51, 1).Value = 1
$ws.Cells.Item(51, 2).Value = 41
$ws.Cells.Item(51, 3).Value = 43
$ws.Cells.Item(51, 4).Value = 170
$ws.Cells.Item(51, 5).Value = 6
$ws.Range("E51").Interior.Color = 5296274
$ws.Cells.Item(51, 6).Value = 0.033

$ws.Cells.Item(52, 1).Value = 2
$ws.Cells.Item(52, 2).Value = 1450
$ws.Cells.Item(52, 3).Value = 1452
$ws.Cells.Item(52, 4).Value = 13303
$ws.Cells.Item(52, 5).Value = 9
$ws.Range("E52").Interior.Color = 5296274
$ws.Cells.Item(52, 6).Value = 3.82

$ws.Cells.Item(53, 1).Value = 3
$ws.Cells.Item(53, 2).Value = 5040
$ws.Cells.Item(53, 3).Value = 5042
$ws.Cells.Item(53, 4).Value = 44944
$ws.Cells.Item(53, 5).Value = 12
$ws.Range("E53").Interior.Color = 5296274
$ws.Cells.Item(53, 6).Value = 16.35

# --- 10. ---------------------------------------------------------------
$ws.Range("A55").NumberFormat = "@"
$ws.Range("A55").Value = "10."
$ws.Range("A55").Style = "Stand."
Set-Header 56

$ws.Cells.Item(57, 1).Value = 1
$ws.Cells.Item(58, 1).Value = 2
$ws.Cells.Item(59, 1).Value = 3

# ---------------------------------------------------------------------
# 4. Restore the active selection to match the author's last position.
# ---------------------------------------------------------------------
$ws.Range("D45").Select()
